# Corrections to the Euclid-proposition worksheet (typo / wording fixes,
# a justification fill-in, a step-numbering fix, row-height reset and the
# resulting selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wording fixes -----------------------------------------------------

# Step 7 (row 8): "el circulo" -> "la circunferencia", "donde se interseca"
# -> "es donde se interseca", and add a trailing period.
$ws.Range("B8").Value = 'Tomar un punto donde $E$ es donde se interseca la circunferencia con la recta.'

# Step 8 (row 9): same fix, for point F.
$ws.Range("B9").Value = 'Tomar un punto donde $F$ es donde se interseca la circunferencia con la recta.'

# Step 13 (row 14): fill in the missing justification.
$ws.Range("C14").Value = 'Postulado 1'

# Steps 19/20 (rows 20-21): accent fix "esta sobre" -> "está sobre".
$ws.Range("B20").Value = 'Como $\overline{GC}$ está sobre $\overline{AB}$ y sus ángulos adyacentes son iguales $\overline{GC}$ es perpendicular a $\overline{AB}$'
$ws.Range("B21").Value = 'Como $\overline{GC}$ está sobre $\overline{AB}$ y sus ángulos adyacentes son iguales $\overline{GC}$ es perpendicular a $\overline{AB}$'

# Step 18 (row 19): triangle notation corrected to angle notation.
$ws.Range("B19").Value = 'Por los pasos (16) y (17) decimos que $\angle EGC = \angle FGC$ '

# --- Step-numbering fix -------------------------------------------------

# Row 21 repeated step number 20 by mistake; it should read 19 (matching
# row 20 directly above it).
$ws.Range("A21").Value = 19

# --- Row height reset ----------------------------------------------------

# Row 19 no longer needs the taller wrapped-text height now that its text
# moved to a single logical line worth of wrapping; reset to the sheet's
# default (auto) height like rows 17/18.
$ws.Rows.Item(19).AutoFit()

# --- Final selection -----------------------------------------------------

$ws.Range("B19").Select()
